$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data row was inserted at row 820 (Uva "Crimpson Seedless" entry),
# pushing all existing rows from 820 downward down by one (last row becomes 926).
$ws.Rows.Item(820).Insert()

# Populate the newly inserted row 820 with its values.
$ws.Cells.Item(820, 1).Value = 9
$ws.Cells.Item(820, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(820, 3).Value = "Metropolitana"
$ws.Cells.Item(820, 4).Value = 45142
$ws.Cells.Item(820, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(820, 5).Value = 13
$ws.Cells.Item(820, 6).Value = "Fruta"
$ws.Cells.Item(820, 7).Value = 100109
$ws.Cells.Item(820, 8).Value = "Uva"
$ws.Cells.Item(820, 9).Value = 100109001
$ws.Cells.Item(820, 10).Value = "Uva"
$ws.Cells.Item(820, 11).Value = "Crimpson Seedless"
$ws.Cells.Item(820, 12).Value = "Primera"
$ws.Cells.Item(820, 13).Value = 130
$ws.Cells.Item(820, 14).Value = 14000
$ws.Cells.Item(820, 15).Value = 14000
$ws.Cells.Item(820, 16).Value = 14000
$ws.Cells.Item(820, 17).Value = "`$/bandeja 8 kilos"
$ws.Cells.Item(820, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(820, 19).Value = 1750
$ws.Cells.Item(820, 20).Value = 8
